$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: update datetime value in column A ---
$ws.Range("A12").Value2 = 38991.45833333334

# --- Row 209: update open/high/low/close values ---
$ws.Range("C209").Value2 = 9059276168200
$ws.Range("D209").Value2 = 9059276168200
$ws.Range("E209").Value2 = 9059276168200
$ws.Range("F209").Value2 = 9059276168200

# --- Row 210: update open/high/low/close values ---
$ws.Range("C210").Value2 = 9505999258000
$ws.Range("D210").Value2 = 9505999258000
$ws.Range("E210").Value2 = 9505999258000
$ws.Range("F210").Value2 = 9505999258000

# --- Row 212: update open/high/low/close values ---
$ws.Range("C212").Value2 = 10834532376000
$ws.Range("D212").Value2 = 10834532376000
$ws.Range("E212").Value2 = 10834532376000
$ws.Range("F212").Value2 = 10834532376000

# --- Row 213: new row appended (copy formatting from row 212 first) ---
$ws.Range("A212:G212").Copy()
$ws.Range("A213:G213").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A213").Value2 = 45108.41666666666
$ws.Range("B213").Value = "ECONOMICS:TRM2"
$ws.Range("C213").Value2 = 11478416920800
$ws.Range("D213").Value2 = 11478416920800
$ws.Range("E213").Value2 = 11478416920800
$ws.Range("F213").Value2 = 11478416920800
$ws.Range("G213").Value2 = 0

$excel.CutCopyMode = 0
